# Order Successpage Code Updated
#
# The "OCT&NonRX_Stage" sheet gets a new lead-in row (copied from the
# "Otcandnonrx" sheet's first row: "Jusdee 400IU Drops 30ml") and its old
# second row ("Keraglo Men Tablet 10'S") is replaced with
# "Folvite Tablet 45'S". The active sheet/selection moves from
# "Otcandnonrx" to "OCT&NonRX_Stage".

$wb = $excel.ActiveWorkbook

$wsOtc = $wb.Worksheets.Item("Otcandnonrx")
$wsStage = $wb.Worksheets.Item("OCT&NonRX_Stage")

# Insert a new first row in the Stage sheet and copy row 1 (cell + style)
# from the Otcandnonrx sheet into it.
$wsStage.Rows.Item(1).Insert()
$wsOtc.Range("A1:B1").Copy($wsStage.Range("A1:B1"))

# The old "Keraglo Men Tablet 10'S" row is now row 3 - replace its text
# with "Folvite Tablet 45'S" (keeping the existing cell style).
$wsStage.Range("A3").Value = "Folvite Tablet 45'S"
$wsStage.Range("B3").Value = "Folvite Tablet 45'S"

# Update active sheet / selections: Otcandnonrx is no longer active, and
# its selection becomes the whole used range (A1:B5); OCT&NonRX_Stage
# becomes the active sheet with C13 selected.
$wsOtc.Activate()
$wsOtc.Range("A1:B5").Select()

$wsStage.Activate()
$wsStage.Range("C13").Select()
